$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header
# formatting (bold/centered/bordered style already used by B1:H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells I2 and J2 (unstyled, like the rest of row 2)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
